# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh values to the Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2747.762
$ws.Range("I2").Value = 5166.8
$ws.Range("J2").Value = 548.63635
$ws.Range("K2").Value = 5166.8
$ws.Range("L2").Value = 548.63635
$ws.Range("M2").Value = -5053.8
$ws.Range("N2").Value = -774.63635
$ws.Range("H18").Value = 858.6
$ws.Range("I18").Value = 858.6
$ws.Range("K18").Value = 858.6
$ws.Range("M18").Value = -574.6
$ws.Range("H33").Value = 394.55554
$ws.Range("I33").Value = 461.16666
$ws.Range("J33").Value = 261.33334
$ws.Range("K33").Value = 461.16666
$ws.Range("L33").Value = 261.33334
$ws.Range("M33").Value = -232.16666
$ws.Range("N33").Value = -719.33334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 875
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H32").Value = 9223.638999999999
$ws.Range("I32").Value = 7566.7812
$ws.Range("J32").Value = 22478.5
$ws.Range("K32").Value = 7566.7812
$ws.Range("L32").Value = 22478.5
$ws.Range("M32").Value = -7279.7812
$ws.Range("N32").Value = -23052.5
$ws.Range("H61").Value = 2434.9443
$ws.Range("J61").Value = 2007
$ws.Range("L61").Value = 2007
$ws.Range("N61").Value = -2431
$ws.Range("H88").Value = 1556.7858
$ws.Range("I88").Value = 1820
$ws.Range("J88").Value = 1410.5555
$ws.Range("K88").Value = 1820
$ws.Range("L88").Value = 1410.5555
$ws.Range("M88").Value = -1414
$ws.Range("N88").Value = -2222.5555
$ws.Range("H91").Value = 1556.7858
$ws.Range("I91").Value = 1820
$ws.Range("J91").Value = 1410.5555
$ws.Range("K91").Value = 1820
$ws.Range("L91").Value = 1410.5555
$ws.Range("M91").Value = -416
$ws.Range("N91").Value = -4218.5555
$ws.Range("H132").Value = 8165.1113
$ws.Range("I132").Value = 8069.4287
$ws.Range("K132").Value = 24208.2861
$ws.Range("M132").Value = -21678.2861
$ws.Range("H136").Value = 2434.9443
$ws.Range("J136").Value = 2007
$ws.Range("L136").Value = 6021
$ws.Range("N136").Value = -11121

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7999.6665
$ws.Range("I20").Value = 7999.6665
$ws.Range("K20").Value = 7999.6665
$ws.Range("M20").Value = -7752.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1002
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 1002
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 1002
$ws.Range("N6").Value = -1228
$ws.Range("M6").ClearContents()
$ws.Range("I22").Value = 1683.75
$ws.Range("J22").Value = 4599.8
$ws.Range("K22").Value = 1683.75
$ws.Range("L22").Value = 4599.8
$ws.Range("M22").Value = -1333.75
$ws.Range("N22").Value = -5299.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1938140
$ws.Range("I4").Value = 4250119.5
$ws.Range("K4").Value = 12750358.5
$ws.Range("M4").Value = -12750246.5
$ws.Range("H34").Value = 2186.5
$ws.Range("J34").Value = 2882.2222
$ws.Range("L34").Value = 8646.6666
$ws.Range("N34").Value = -8814.6666
$ws.Range("H109").Value = 145550.28
$ws.Range("I109").Value = 202135.4
$ws.Range("J109").Value = 4087.5
$ws.Range("K109").Value = 606406.2
$ws.Range("L109").Value = 12262.5
$ws.Range("M109").Value = -605366.2
$ws.Range("N109").Value = -14342.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 3555725
$ws.Range("I14").Value = 5333433.5
$ws.Range("J14").Value = 308.33334
$ws.Range("K14").Value = 5333433.5
$ws.Range("L14").Value = 308.33334
$ws.Range("M14").Value = -5333265.5
$ws.Range("N14").Value = -644.33334
$ws.Range("H80").Value = 5752.5
$ws.Range("I80").Value = 1499
$ws.Range("J80").Value = 10006
$ws.Range("K80").Value = 1499
$ws.Range("L80").Value = 10006
$ws.Range("M80").Value = -501
$ws.Range("N80").Value = -12002
$ws.Range("H83").Value = 5752.5
$ws.Range("I83").Value = 1499
$ws.Range("J83").Value = 10006
$ws.Range("K83").Value = 7495
$ws.Range("L83").Value = 50030
$ws.Range("M83").Value = -2503
$ws.Range("N83").Value = -60014
$ws.Range("H126").Value = 5167.737
$ws.Range("I126").Value = 3924.4167
$ws.Range("J126").Value = 7299.143
$ws.Range("K126").Value = 11773.2501
$ws.Range("L126").Value = 21897.429
$ws.Range("M126").Value = -9303.250100000001
$ws.Range("N126").Value = -26837.429
$ws.Range("H132").Value = 75594.31
$ws.Range("I132").Value = 116551.1
$ws.Range("K132").Value = 349653.3
$ws.Range("M132").Value = -347123.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6576.091
$ws.Range("I7").Value = 5460.6
$ws.Range("J7").Value = 7505.6665
$ws.Range("K7").Value = 5460.6
$ws.Range("L7").Value = 7505.6665
$ws.Range("M7").Value = -5348.6
$ws.Range("N7").Value = -7729.6665
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H55").Value = 1327.4783
$ws.Range("I55").Value = 1058.8
$ws.Range("J55").Value = 1534.1538
$ws.Range("K55").Value = 1058.8
$ws.Range("L55").Value = 1534.1538
$ws.Range("M55").Value = -885.8
$ws.Range("N55").Value = -1880.1538
$ws.Range("H82").Value = 4262.4546
$ws.Range("I82").Value = 2213.6667
$ws.Range("J82").Value = 5030.75
$ws.Range("K82").Value = 2213.6667
$ws.Range("L82").Value = 5030.75
$ws.Range("M82").Value = -1852.6667
$ws.Range("N82").Value = -5752.75
$ws.Range("H85").Value = 4262.4546
$ws.Range("I85").Value = 2213.6667
$ws.Range("J85").Value = 5030.75
$ws.Range("K85").Value = 2213.6667
$ws.Range("L85").Value = 5030.75
$ws.Range("M85").Value = -965.6667000000002
$ws.Range("N85").Value = -7526.75
$ws.Range("H93").Value = 1741.25
$ws.Range("J93").Value = 1806
$ws.Range("L93").Value = 1806
$ws.Range("N93").Value = -4302
$ws.Range("H126").Value = 6576.091
$ws.Range("I126").Value = 5460.6
$ws.Range("J126").Value = 7505.6665
$ws.Range("K126").Value = 16381.8
$ws.Range("L126").Value = 22516.9995
$ws.Range("M126").Value = -13911.8
$ws.Range("N126").Value = -27456.9995
$ws.Range("H132").Value = 4743.125
$ws.Range("I132").Value = 1725
$ws.Range("J132").Value = 5749.1665
$ws.Range("K132").Value = 5175
$ws.Range("L132").Value = 17247.4995
$ws.Range("M132").Value = -2645
$ws.Range("N132").Value = -22307.4995
$ws.Range("H139").Value = 76000
$ws.Range("J139").Value = 105000
$ws.Range("L139").Value = 105000
$ws.Range("N139").Value = -115280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 16675000
$ws.Range("I5").Value = 30000000
$ws.Range("J5").Value = 6015000.5
$ws.Range("K5").Value = 30000000
$ws.Range("L5").Value = 6015000.5
$ws.Range("M5").Value = -29999888
$ws.Range("N5").Value = -6015224.5
